$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 196
$ws.Range("B17").Value = 'Lots of people with high profiles tweeting stuff in the last 12 hours about Australia and Australians they will end up regretting...#ausvotes'
$ws.Range("C17").Value = -0.75825214
$ws.Range("D17").Value = 'Melbourne'

$ws.Range("A18").Value = 190
$ws.Range("B18").Value = 'Newly unemployed Abbott wondering which idiot screwed up Centrelink this badly https://t.co/PxWK4KqCX2 #auspol https://t.co/uCaSvl8ZCs'
$ws.Range("C18").Value = -0.8318004
$ws.Range("D18").Value = 'Sydney'

$ws.Range("A19").Value = 174
$ws.Range("B19").Value = 'Things aren''t looking very good at Labor HQ. Here, the mood has shifted quite dramatically. LIVE updates: https://t.co/ZL03xBTyyM #AusVotes #Auspol https://t.co/btzOwZsfk2'
$ws.Range("C19").Value = -0.8639997
$ws.Range("D19").Value = 'Australia'

$ws.Range("A20").Value = 172
$ws.Range("B20").Value = 'I would rather live next door to Muslims any day.. I''ve already lived next door to some bogans and it was bloody terrible ...drinking fighting and loud music.. Anning is evil    #auspol #AusVotes2019 https://t.co/DuzG4WfVGo'
$ws.Range("C20").Value = -0.93293774
$ws.Range("D20").Value = 'brisbane'

$ws.Range("A21").Value = 161
$ws.Range("B21").Value = 'Trump politics where political parties &amp; billionaires influence elections by spreading lies via social media to influence elections has reached Australia in a big way. I am being constantly told by people they have received fake news about Labor''s tax policies #ausvotes'
$ws.Range("C21").Value = -0.8363542
$ws.Range("D21").Value = 'Melbourne, Victoria'

$ws.Range("A22").Value = 145
$ws.Range("B22").Value = '.@TonyAbbottMHR has lost his seat!!!!!
Good. Bloody. Riddance. 
#ausvotes #auspol #TonyAbbott'
$ws.Range("C22").Value = -0.8743961
$ws.Range("D22").Value = 'Australia, Wurundjeri Land'

$ws.Range("A23").Value = 142
$ws.Range("B23").Value = 'Fmr PM @TonyAbbottMHR says the campaign in #Warringah has been ''pretty ugly'' and he''s urged voters not to give "these nasty elements... the victory they want" #ausvotes https://t.co/NxItnJYVkZ'
$ws.Range("C23").Value = -0.7992251500000001
$ws.Range("D23").Value = 'Australia'

$ws.Range("A24").Value = 142
$ws.Range("B24").Value = 'Is there any man in Australia who deserves Chris Bowen''s fate tonight
Alan Jones is going off like a cranky macaw. 
#Auspol #Ausvotes2019 #Channel7 https://t.co/xsVaDX6OZx'
$ws.Range("C24").Value = -0.8074239
$ws.Range("D24").Value = 'Australia'

$ws.Range("A25").Value = 141
$ws.Range("B25").Value = 'If you think #Barnaby Christensen &amp; Angus Taylor were crook before this...you watch em go now. They just got the green light they wanted and will loot us all for everything. #auspol'
$ws.Range("C25").Value = -0.85240823
$ws.Range("D25").Value = 'The Adelaide Diaspora'

$ws.Range("A26").Value = 138
$ws.Range("B26").Value = 'Today I’m sad for so many reasons, but it is the thought that last night sealed the fate and ruined the future of two innocent Australian born children that really breaks my heart. Australia, I thought you were better than this... #auspol https://t.co/pqoGqc7DwZ'
$ws.Range("C26").Value = -0.9233246000000001
$ws.Range("D26").Value = 'Brisbane, Australia'

$ws.Range("A27").Value = 137
$ws.Range("B27").Value = 'This election was won on lies, hollow promises, pork-barreling &amp; slogans that are impossible to deliver upon.
I don''t think we''ll have to wait 3 yrs to go back to the polls, especially with a hostile Senate
#auspol'
$ws.Range("C27").Value = -0.8365116
$ws.Range("D27").Value = 'Australia'

$ws.Range("A28").Value = 135
$ws.Range("B28").Value = 'BREAKING NEWS. I''ve been engaging with @GatekeepersReef for 11 months. They are FAKE and seemingly linked to the Libs. Please let people know by sharing this expose by @jrojourno for @abcnews. Then report the page! #StopAdani #auspol #AusVotes2019 
https://t.co/8VpIpf5v3J'
$ws.Range("C28").Value = -0.86193025
$ws.Range("D28").Value = 'Brisbane, Queensland'

$ws.Range("A29").Value = 126
$ws.Range("B29").Value = 'So @PeterDutton_MP was happy to import axe-murderers, and let them roam free among us while, simultaneously, accusing Labor of encouraging rapists &amp; terrorists to risk a boat to Oz?
#Dicksonvotes needs to allow us to prosecute this man.
#ausvotes
https://t.co/LB6X7TmkG0'
$ws.Range("C29").Value = -0.8326862
$ws.Range("D29").Value = '#Tinonee, NSW, Oz'

$ws.Range("A30").Value = 125
$ws.Range("B30").Value = 'Leigh Sales: "If Newspoll has it so wrong, isn''t it time you stop sacking your leaders over it?" #ausvotes'
$ws.Range("C30").Value = -0.768327
$ws.Range("D30").Value = 'Sydney, Australia'

$ws.Range("A31").Value = 125
$ws.Range("B31").Value = 'The reason the Liberals haven''t rolled out any policies during this election campaign is because the warring factions within the party can''t agree on anything that should be policy. A vote for the LNP is a vote for another 3 years of absolute chaos. #auspol'
$ws.Range("C31").Value = -0.8473929
$ws.Range("D31").Value = 'Brisbane, Queensland'

$ws.Range("A32").Value = 115
$ws.Range("B32").Value = 'I officially hate anyone who voted for the Australian Liberals/Nationals in this year''s election. And yes, I would allow politics to get in the way of family and friendship. You essentially voted to take away the humanity of this country''s most vulnerable &amp; marginalised. #AusPol'
$ws.Range("C32").Value = -0.93352586
$ws.Range("D32").Value = 'Melbourne, Victoria'

$ws.Range("A33").Value = 115
$ws.Range("B33").Value = 'The Happy Clapper Slogan Bogan Liar from the Shire says he will cut “Green Tape”? Do LNP voters have any idea what that means?! Well they bloody well should! They are environment protection laws. Care about global warming? Why the hell would you vote LNP? Are they schizo? #auspol'
$ws.Range("C33").Value = -0.8778286
$ws.Range("D33").Value = 'Sydney, New South Wales'

$ws.Range("A34").Value = 113
$ws.Range("B34").Value = '"You’re not children any more. I didn’t mind explaining photosynthesis to you when you were 12. But you’re adults now, and this is an actual crisis! Got it?! Safety glasses off,  motherfuckers" 
TY @BillNye for channeling how we feel here in Aus #ausvotes
https://t.co/BYxjKNt3Ro'
$ws.Range("C34").Value = -0.8491100700000001
$ws.Range("D34").Value = 'Wurundjeri Country / Melbourne'

$ws.Range("A35").Value = 101
$ws.Range("B35").Value = 'A vote for the Liberals is a vote for more chaos. #auspol https://t.co/rTqk8D7UkL'
$ws.Range("C35").Value = -0.76383
$ws.Range("D35").Value = 'Australia'

$ws.Range("A36").Value = 99
$ws.Range("B36").Value = 'journalists at the Australian are uncomfortable with the paper’s cheerleading for the Coalition. Wow it''s worse at the Courier-Mail I''ve been told from someone, it''s almost at boiling point. Almost ready for a walkout.  Not everyone who works there votes liberal  ffs #auspol'
$ws.Range("C36").Value = -0.900532
$ws.Range("D36").Value = 'brisbane'

$ws.Range("A37").Value = 97
$ws.Range("B37").Value = 'This is what you can expect if Scott Morrison’s coalition of chaos wins. #auspol https://t.co/qt5SZki2Py'
$ws.Range("C37").Value = -0.8277828
$ws.Range("D37").Value = 'Australia'

$ws.Range("A38").Value = 97
$ws.Range("B38").Value = 'These unAustralian idiots just don''t get it. They hate Australia and Australians and probably themselves too. This is why #Labor lost. And this is why #Labor will always lose while it panders to this sorry lot of haters. 
#AusVotes2019 #auspol #AustraliaDecides #ausvotes https://t.co/NqbwHEb2Ud'
$ws.Range("C38").Value = -0.9304413
$ws.Range("D38").Value = 'Victoria, Australia'

$ws.Range("A39").Value = 92
$ws.Range("B39").Value = 'This is so true, Labor will Tax the people of Australia into economic oblivion. #auspol #AusVotes19 #insiders #LaborTrash 
https://t.co/h1gYgy6PJj'
$ws.Range("C39").Value = -0.87745017
$ws.Range("D39").Value = 'Australia'

$ws.Range("A40").Value = 89
$ws.Range("B40").Value = 'I find the height of rudeness is when someone puts their hand out to shake another’s &amp; say Gday, &amp; they see you are the Labor Candidate &amp; say ‘nope’! I mean I don’t mind if you have a different view but to not shake a hand is just disrespectful. You can still be civil. #auspol'
$ws.Range("C40").Value = -0.74149895
$ws.Range("D40").Value = 'Queensland, Australia'

$ws.Range("A41").Value = 88
$ws.Range("B41").Value = 'Reasons not to vote for Scott Morrison: 
-Voted against Banking Royal Commission 26 times
-Abstained from the same sex marriage vote 
-Voted to cut penalty rates 
- Urged the shadow cabinet to capitalise on “growing concerns of Muslims”
I could go all day...
#auspol #ausvotes'
$ws.Range("C41").Value = -0.7185315
$ws.Range("D41").Value = 'Canberra, Australian Capital Territory'

$ws.Range("A42").Value = 86
$ws.Range("B42").Value = 'UNHINGED hysteria. 
WITHOUT cheap, abundant hydrocarbons aka #fossilfuels, every tree on the planet would have been cut down by now for heating, cooking and building needs. 
#ClimateChange #Delusions #Misanthropy #Auspol #AusVotes2019 #WarringahVotes #CDNpoli https://t.co/WZIXpnuyYF'
$ws.Range("C42").Value = -0.7570373
$ws.Range("D42").Value = 'Melbourne, Victoria'

$ws.Range("A43").Value = 85
$ws.Range("B43").Value = 'Thousands of workers have had their take home pay slashed under the Liberals. With further cuts to penalty rates due to hit in coming years, can Australians really afford three more years of Liberal Chaos? #auspol https://t.co/eaB4CVx9ob'
$ws.Range("C43").Value = -0.8845269
$ws.Range("D43").Value = 'Australia'

$ws.Range("A44").Value = 83
$ws.Range("B44").Value = 'This misleading post on WeChat has been circulating widely today, falsely claiming Labor will implement a high level inheritance tax. Labor has already complained to FB about similar claims on that platform. Translation thanks to @ABCChinese @abcnews #ausvotes #ausvotes19 #auspol https://t.co/s4DWktGIC6'
$ws.Range("C44").Value = -0.82743406
$ws.Range("D44").Value = 'Sydney, Australia'

$ws.Range("A45").Value = 80
$ws.Range("B45").Value = 'I''m sitting here on Twitter  on the lounge not game enough to turn on a radio the TV or read online newspapers my world just got a lot fucking smaller   #auspol going for a bike ride soon 😘😢'
$ws.Range("C45").Value = -0.84640074
$ws.Range("D45").Value = 'brisbane'

$ws.Range("A46").Value = 79
$ws.Range("B46").Value = 'Can’t afford to pay his workers, can afford to piss off to Fiji in the middle of an election campaign. Clive Palmer is a joke, and every person conned into voting for him is the punchline. #Auspol #AusVotes19 https://t.co/ojbh9Fp97Z'
$ws.Range("C46").Value = -0.9174092
$ws.Range("D46").Value = 'Ballarat'

$ws.Range("A47").Value = 78
$ws.Range("B47").Value = '#zalies
Zali says she wants to run a respectful campaign, but instead, her campaigners are engaging in these sorts of grubby, juvenile and dirty tactics and actions
But Zali the notorious hypocrite won''t call them out! Imagine if this was a Tony backer!
#warringahvotes #auspol https://t.co/sxdM5XKOoo'
$ws.Range("C47").Value = -0.82532406
$ws.Range("D47").Value = 'Sydney, New South Wales'

$ws.Range("A48").Value = 77
$ws.Range("B48").Value = 'The Liberals have caused massive wage stagnation, thousands in precarious employment, they have overseen abhorrent discrimination within CDEP program, $4/Hr Interns under PATH, attacked single parents and seen suicides over #robodebt &amp; ALL THEY CARE ABOUT IS THE WEALTHY! #Auspol https://t.co/LBjWxxYFck'
$ws.Range("C48").Value = -0.9424369
$ws.Range("D48").Value = 'Rockhampton'

$ws.Range("A49").Value = 72
$ws.Range("B49").Value = 'This has also been a story doing the rounds for eons. Turnbull’s office covered this one up.
STABLE GOVT? My arse!               
#AusVotes2019 #Auspol 
“Two of MP''s former staffers say they asked Scott Morrison to help resolve their complaints https://t.co/fupuBi43hK'
$ws.Range("C49").Value = -0.83152956
$ws.Range("D49").Value = 'Sydney'

$ws.Range("A50").Value = 71
$ws.Range("B50").Value = '#BobHawkeRIP  I''m starting to think nobody wants to talk about Hawke Government reduction in poverty 80% in couples without work and kids 
50% for single-parent family''s Nobody has done this before or since 
Why will we not mention this #auspol'
$ws.Range("C50").Value = -0.7662709
$ws.Range("D50").Value = 'Geelong, Victoria'

$ws.Range("A51").Value = 71
$ws.Range("B51").Value = 'Australians all let us revoice
Coz it’s plain to see
We’re old and soiled in mortal coils
And we suffer the politics of sociopathy
Our land abounds with racist gits
Of religious psychopaths
In history''s page maintain the rage 
Let’s make Australia fair...
#ausvotes https://t.co/RdEXPRYGfc'
$ws.Range("C51").Value = -0.72214097
$ws.Range("D51").Value = '#Tinonee, NSW, Oz'
